$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title / header text updates (in-place shared-string edits) ----
# "Volume 30   Number  25" -> "...  26"  (digits at characters 21-22)
$ws.Range("A8").Characters(21, 2).Text = "26"

# "Report Covering the Week  6/19/2023  Through  6/25/2023"
#  -> "...6/26/2023  Through  7/2/2023"
$ws.Range("C9").Characters(27, 9).Text = "6/26/2023"
$ws.Range("C9").Characters(47, 9).Text = "7/2/2023"

# ---- Reference cells used to clone formatting (style) onto cells whose
#      value "type" (number vs text placeholder) changes this week ----
$styleRefText = $ws.Range("A14")
$styleRefNum1 = $ws.Range("F14")
$styleRefNum2 = $ws.Range("K14")

# ---- Simple numeric value updates (no type/style change) ----
$ws.Range("N14").Value = -60
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 22
$ws.Range("K15").Value = 37.5
$ws.Range("L15").Value = 46.666666666666
$ws.Range("M15").Value = 22.222222222222
$ws.Range("N15").Value = -53.191489361702
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -4
$ws.Range("I16").Value = 133
$ws.Range("J16").Value = 176
$ws.Range("K16").Value = -24.431818181818
$ws.Range("L16").Value = 24.29906542056
$ws.Range("M16").Value = -34.482758620689
$ws.Range("N16").Value = -89.044481054365
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -71.428571428571
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = -19.642857142857
$ws.Range("I17").Value = 325
$ws.Range("J17").Value = 308
$ws.Range("K17").Value = 5.51948051948
$ws.Range("L17").Value = 7.97342192691
$ws.Range("M17").Value = 47.058823529411
$ws.Range("N17").Value = -47.325769854132
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 107
$ws.Range("K18").Value = 0.934579439252
$ws.Range("L18").Value = 63.636363636363
$ws.Range("M18").Value = -41.935483870967
$ws.Range("N18").Value = -89.605389797882
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 8.333333333333
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 21.56862745098
$ws.Range("I19").Value = 349
$ws.Range("J19").Value = 353
$ws.Range("K19").Value = -1.13314447592
$ws.Range("L19").Value = 78.974358974359
$ws.Range("M19").Value = 30.223880597014
$ws.Range("N19").Value = -6.684491978609
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = -35.483870967741
$ws.Range("I20").Value = 128
$ws.Range("J20").Value = 148
$ws.Range("K20").Value = -13.513513513513
$ws.Range("L20").Value = 20.754716981132
$ws.Range("M20").Value = -10.48951048951
$ws.Range("N20").Value = -87.981220657277
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -32
$ws.Range("G21").Value = 184
$ws.Range("H21").Value = -4.891304347826
$ws.Range("I21").Value = 1073
$ws.Range("J21").Value = 1117
$ws.Range("K21").Value = -3.939122649955
$ws.Range("L21").Value = 34.798994974874
$ws.Range("M21").Value = 2.385496183206
$ws.Range("N21").Value = -75.479890310786
$ws.Range("F23").Value = 1
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = -33.333333333333
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 12
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -2.298850574712
$ws.Range("I24").Value = 612
$ws.Range("J24").Value = 488
$ws.Range("K24").Value = 25.409836065573
$ws.Range("L24").Value = 36.607142857142
$ws.Range("M24").Value = 29.113924050632
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -31.818181818181
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = -23.456790123456
$ws.Range("I25").Value = 381
$ws.Range("J25").Value = 376
$ws.Range("K25").Value = 1.329787234042
$ws.Range("L25").Value = 23.701298701298
$ws.Range("M25").Value = -17.353579175705
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 7
$ws.Range("H26").Value = 250
$ws.Range("I26").Value = 28
$ws.Range("K26").Value = 21.739130434782
$ws.Range("L26").Value = 16.666666666666
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 38
$ws.Range("K27").Value = -7.317073170731
$ws.Range("L27").Value = 15.151515151515
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -28.125
$ws.Range("M28").Value = -23.333333333333
$ws.Range("N28").Value = -76.288659793814
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -32.142857142857
$ws.Range("M29").Value = -26.923076923076
$ws.Range("N29").Value = -79.347826086956

# ---- Number -> Text placeholder conversions (style 14) ----
$dst = $ws.Range("G14")
$dst.NumberFormat = "@"
$dst.Value = "0"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("H14")
$dst.NumberFormat = "@"
$dst.Value = "***.*"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("C22")
$dst.NumberFormat = "@"
$dst.Value = "0"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("C23")
$dst.NumberFormat = "@"
$dst.Value = "0"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("D27")
$dst.NumberFormat = "@"
$dst.Value = "0"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("E27")
$dst.NumberFormat = "@"
$dst.Value = "***.*"
$styleRefText.Copy()
$dst.PasteSpecial(-4122)

# ---- Text placeholder -> Number conversions ----
$dst = $ws.Range("D23")
$dst.Value = 1
$styleRefNum1.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("E23")
$dst.Value = -100
$styleRefNum2.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("G23")
$dst.Value = 1
$styleRefNum1.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("H23")
$dst.Value = 0
$styleRefNum2.Copy()
$dst.PasteSpecial(-4122)
$dst = $ws.Range("M23")
$dst.Value = 300
$styleRefNum2.Copy()
$dst.PasteSpecial(-4122)

$app.CutCopyMode = $false
